$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "password" column values (G2:G6) ---
# The old demo password "Mayorista2021*.*" is replaced everywhere by "testPassword".
$ws.Range("G2:G6").Value = "testPassword"

# --- Turn the A1:G6 range into a proper Excel Table ("Tabla1") ---
$range = $ws.Range("A1:G6")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Tabla1"
$tbl.TableStyle = "TableStyleLight8"

# Keep the "user" column's data cells on the hyperlink cell style (matches
# the pre-existing mailto: hyperlink formatting used on those cells).
$tbl.ListColumns.Item(6).DataBodyRange.Style = "Hipervínculo"

# --- Resize the columns to fit the table's (wider) content ---
$ws.Columns.Item(1).ColumnWidth = 8.333333333333334
$ws.Columns.Item(2).ColumnWidth = 26
$ws.Columns.Item(4).ColumnWidth = 38
$ws.Columns.Item(5).ColumnWidth = 10.166666666666666
$ws.Columns.Item(6).ColumnWidth = 24.666666666666668
$ws.Columns.Item(7).ColumnWidth = 15

# --- Move the active selection ---
$ws.Range("G7").Select() | Out-Null
